$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-10, excluding row 9) got re-dated / re-priced.
# This corresponds to a cyclic re-shuffle of the D (Fecha), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# R (Origen) and S (Precio $/Kg) values between rows:
#   row 2 <- old row 7
#   row 3 <- old row 10
#   row 4 <- old row 5
#   row 5 <- old row 2
#   row 6 <- old row 4
#   row 7 <- old row 3
#   row 8 <- old row 6
#   row 10 <- old row 8

$rows = @{
    2  = @{ D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 }
    3  = @{ D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 }
    4  = @{ D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 }
    5  = @{ D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";     S = 5000 }
    6  = @{ D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
    7  = @{ D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";     S = 3000 }
    8  = @{ D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";     S = 4500 }
    10 = @{ D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($r, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($r, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($r, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($r, 19).Value = $vals.S   # S: Precio $/Kg
}
